# Add data for 2022-06-22: bumps the "through June NN" running-month
# label/sheet name by one day, and increments the June-column carjacking
# counts for the neighborhoods/years that received a new incident.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell (column B, row 1) and sheet tab name both track the
# "through" date, which advances from June 13 to June 14.
$ws.Range("B1").Value = "June 2022 (through June 14)"
$ws.Name = "Through 2022-06-14"

# Per-neighborhood / per-year "June" counts that changed.
$ws.Range("B3").Value = 5     # Auburn Gresham   - June 2022 (through June 14)
$ws.Range("B4").Value = 5     # North Lawndale   - June 2022 (through June 14)
$ws.Range("H4").Value = 2     # North Lawndale   - June 2021
$ws.Range("N6").Value = 2     # Humboldt Park    - June 2020
$ws.Range("H9").Value = 2     # Grand Crossing   - June 2021
$ws.Range("AR10").Value = 1   # Garfield Park    - June 2015
$ws.Range("N12").Value = 2    # Roseland         - June 2020
$ws.Range("T12").Value = 2    # Roseland         - June 2019
$ws.Range("H14").Value = 6    # Austin           - June 2021
$ws.Range("N14").Value = 4    # Austin           - June 2020
$ws.Range("T14").Value = 1    # Austin           - June 2019
$ws.Range("Z14").Value = 3    # Austin           - June 2018
$ws.Range("AF22").Value = 1   # Chinatown        - June 2017
$ws.Range("T23").Value = 1    # Little Village   - June 2019
$ws.Range("AR23").Value = 1   # Little Village   - June 2015
$ws.Range("N29").Value = 2    # West Ridge       - June 2020
$ws.Range("H33").Value = 1    # Near South Side  - June 2021
$ws.Range("B37").Value = 2    # Fuller Park      - June 2022 (through June 14)
$ws.Range("B40").Value = 2    # Bridgeport       - June 2022 (through June 14)
$ws.Range("B70").Value = 3    # Loop             - June 2022 (through June 14)
$ws.Range("Z75").Value = 1    # New City         - June 2018
$ws.Range("AL95").Value = 1   # West Town        - June 2016
